$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A73").Value = "2023-12-07 15:13:39"
$ws.Range("B73").Value = 0.0008

$ws.Range("A74").Value = "2023-12-07 15:13:52"
$ws.Range("B74").Value = 0.0002
